# card_export.xlsx edit script
# - Adds a new column I (value 0) to rows 1-51 (stat column that already
#   exists for rows 52-79, but was missing for the first block of cards).
# - Adds a "Magie" type (new shared string) in column D plus zeroed stat
#   columns F:I (and E for row 78, which previously had no stats at all)
#   to the spell rows 69-79.
# - Updates the sheet view's scroll position / selection to I1:I51.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ---------------------------------------------------------------------
# 1) New column I for rows 1-51, matching the formatting already used
#    by column H on each respective row, value 0.
# ---------------------------------------------------------------------
$ws.Range("H1").Copy()
$ws.Range("I1:I51").PasteSpecial(-4122)   # xlPasteFormats

for ($r = 1; $r -le 51; $r++) {
    $ws.Cells.Item($r, 9).Value = 0
}

# ---------------------------------------------------------------------
# 2) Spell rows 69-79: add a "Magie" type in column D and zero out the
#    stat columns F, G, H, I (plus E for row 78 which had none yet).
# ---------------------------------------------------------------------
$ws.Range("E69").Copy()
$ws.Range("D69:D79").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E78").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("F69:I79").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D69:D79").Value = "Magie"

for ($r = 69; $r -le 79; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
    $ws.Cells.Item($r, 9).Value = 0
}
$ws.Range("E78").Value = 0

# ---------------------------------------------------------------------
# 3) Update view: scroll near the bottom block and select the newly
#    added column I values.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I1:I51").Select()

$wb.Save()
